$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.595.67"
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").Value = "2.166.86"
$ws.Range("E3").Value = "  +0.36%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'226.60"
$ws.Range("E5").Value = "  -1.38%  "
$ws.Range("D6").Value = "'0.622"
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").Value = "'62.97"
$ws.Range("E7").Value = "  -0.34%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("E9").Value = "  -1.10%  "
$ws.Range("E10").Value = "  -0.92%  "
$ws.Range("E11").Value = "  -0.01%  "
$ws.Range("D12").Value = "'15.86"
$ws.Range("E12").Value = "  -1.94%  "
$ws.Range("D13").Value = "2.489.02"
$ws.Range("E13").Value = "  +0.51%  "
$ws.Range("D14").Value = "'21.71"
$ws.Range("E14").Value = "  -2.67%  "
$ws.Range("D15").Value = "'0.809"
$ws.Range("E15").Value = "  -1.28%  "
$ws.Range("E16").Value = "  -1.64%  "
$ws.Range("D17").Value = "2.171.52"
$ws.Range("E17").Value = "  +0.30%  "
$ws.Range("D18").Value = "39.557.08"
$ws.Range("E18").Value = "  +0.29%  "
$ws.Range("D19").Value = "0.0₃0915"
$ws.Range("E19").Value = "  +7.21%  "
$ws.Range("D20").Value = "'71.67"
$ws.Range("E20").Value = "  -1.07%  "
$ws.Range("D21").Value = "'6.01"
$ws.Range("E21").Value = "  -2.70%  "
$ws.Range("D22").Value = "'227.53"
$ws.Range("E22").Value = "  -0.76%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").Value = "'2.33"
$ws.Range("E24").Value = "  -1.85%  "
$ws.Range("D25").Value = "'2.33"
$ws.Range("E25").Value = "  -1.55%  "
$ws.Range("D26").Value = "'170.83"
$ws.Range("E26").Value = "  -1.25%  "
$ws.Range("E27").Value = "  -3.36%  "
$ws.Range("E28").Value = "  +0.13%  "
$ws.Range("E29").Value = "  +0.68%  "
$ws.Range("D30").Value = "'19.67"
$ws.Range("E30").Value = "  +0.12%  "
$ws.Range("D31").Value = "'2.67"
$ws.Range("E31").Value = "  +3.87%  "
$ws.Range("E32").Value = "  -0.15%  "
$ws.Range("E33").Value = "  -3.37%  "
$ws.Range("D34").Value = "'4.70"
$ws.Range("E34").Value = "  -3.11%  "
$ws.Range("E35").Value = "  -3.29%  "
$ws.Range("E36").Value = "  -1.39%  "
$ws.Range("D37").Value = "'3.81"
$ws.Range("E37").Value = "  +6.17%  "
$ws.Range("E38").Value = "  -1.02%  "
$ws.Range("E39").Value = "  -0.10%  "
$ws.Range("D40").Value = "'4.94"
$ws.Range("E40").Value = "  +18.87%  "
$ws.Range("D41").Value = "'102.48"
$ws.Range("E41").Value = "  -1.06%  "
$ws.Range("E42").Value = "  -1.59%  "
$ws.Range("E43").Value = "  -3.70%  "
$ws.Range("D44").Value = "1.513.31"
$ws.Range("E44").Value = "  -1.86%  "
$ws.Range("E45").Value = "  +1.34%  "
$ws.Range("D46").Value = "'7.90"
$ws.Range("E46").Value = "  +1.04%  "
$ws.Range("D47").Value = "'0.0921"
$ws.Range("E47").Value = "  -0.61%  "
$ws.Range("E48").Value = "  -0.17%  "
$ws.Range("E49").Value = "  -1.95%  "
$ws.Range("E50").Value = "  +31.54%  "
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.368.26"
$ws.Range("E51").Value = "  +0.32%  "
